$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.945.58"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.817.08"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.305"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.23%  "
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "2.081.05"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "1.808.57"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("E15").Value = "  +5.76%  "
$ws.Range("E16").Value = "  +5.04%  "
$ws.Range("D17").Value = "34.923.21"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.21%  "
$ws.Range("D19").Value = "0.0₃0786"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.11%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +30.33%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").Value = "3.341.29"
$ws.Range("E31").Value = "  +37.52%  "
$ws.Range("E32").Value = "  +6.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "93.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.675"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.76%  "
$ws.Range("E38").Value = "  +5.46%  "
$ws.Range("D39").Value = "1.308.40"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0191"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.980"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0511"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").Value = "1.992.05"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
